$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.947550356388092
$ws.Range("B1").Value = 1.617266774177551
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.565396308898926
$ws.Range("E1").Value = 1.346478939056396
